# Lista Sendo Atualizada Online e Offline
# Adds newly-onboarded customers to both the "Clientes Belo Horizonte Online"
# and "Clientes Belo Horizonte Offline" sheets, fills in a couple of
# previously-blank "Nicho" cells, and nudges a couple of cosmetic view
# properties (column widths / scroll position) to match where the author
# had scrolled to while editing.

$wb = $excel.ActiveWorkbook
$wsOnline = $wb.Worksheets.Item("Clientes Belo Horizonte Online")
$wsOffline = $wb.Worksheets.Item("Clientes Belo Horizonte Offline")

# ---------------------------------------------------------------------
# Sheet "Clientes Belo Horizonte Online" — new rows 43-49
# ---------------------------------------------------------------------

# Row 43
$wsOnline.Range("A43").Value = 38
$wsOnline.Range("B43").Value = "Sociedade Protetora dos Animais"
$wsOnline.Range("E43").Value = "www.sociedadeprotetoradosanimaisbh.org"
$wsOnline.Hyperlinks.Add($wsOnline.Range("E43"), "http://www.sociedadeprotetoradosanimaisbh.org", "", "", "www.sociedadeprotetoradosanimaisbh.org") | Out-Null
$wsOnline.Range("F43").Value = "Cães"

# Row 44
$wsOnline.Range("A44").Value = 39
$wsOnline.Range("B44").Value = "Adoção Bh"
$wsOnline.Range("D44").Value = "(31) 3081-0833"
$wsOnline.Range("E44").Value = "adoçãobh@gmail.com"
$wsOnline.Hyperlinks.Add($wsOnline.Range("E44"), "mailto:adoçãobh@gmail.com") | Out-Null
$wsOnline.Range("F44").Value = "Cães"

# Row 45
$wsOnline.Range("A45").Value = 40
$wsOnline.Range("B45").Value = "Adote um Amigo"
$wsOnline.Range("D45").Value = "(31) 2535-2517"
$wsOnline.Range("E45").Value = "adote1amigobh@gmail.com"
$wsOnline.Hyperlinks.Add($wsOnline.Range("E45"), "mailto:adote1amigobh@gmail.com") | Out-Null
$wsOnline.Range("F45").Value = "Cães e Gatos"

# Row 46
$wsOnline.Range("A46").Value = 41
$wsOnline.Range("B46").Value = "Asas e Amigos"
$wsOnline.Range("D46").Value = "(31) 9331-2126"
$wsOnline.Range("E46").Value = "asaseamigos@hotmail.com"
$wsOnline.Hyperlinks.Add($wsOnline.Range("E46"), "mailto:asaseamigos@hotmail.com") | Out-Null
$wsOnline.Range("F46").Value = "Animais em geral"

# Row 47
$wsOnline.Range("A47").Value = 42
$wsOnline.Range("B47").Value = "Animaisos Ong"
$wsOnline.Range("E47").Value = "www.animaisos.org"
$wsOnline.Hyperlinks.Add($wsOnline.Range("E47"), "http://www.animaisos.org", "", "", "www.animaisos.org") | Out-Null
$wsOnline.Range("F47").Value = "Animais em geral"

# Row 48
$wsOnline.Range("A48").Value = 43
$wsOnline.Range("B48").Value = "Cão Viver"
$wsOnline.Range("C48").Value = "R. 1º de Maio, nº 165 - Braunas"
$wsOnline.Range("E48").Value = "www.caoviver.com.br"
$wsOnline.Hyperlinks.Add($wsOnline.Range("E48"), "http://www.caoviver.com.br", "", "", "www.caoviver.com.br") | Out-Null
$wsOnline.Range("F48").Value = "Cães e Gatos"

# Row 49
$wsOnline.Range("A49").Value = 44
$wsOnline.Range("B49").Value = "Bast Adotar"

# ---------------------------------------------------------------------
# Sheet "Clientes Belo Horizonte Offline" — fill two "Nicho" gaps and
# add the matching new rows 42-43
# ---------------------------------------------------------------------

$wsOffline.Range("F6").Value = "Veterinário"
$wsOffline.Range("F41").Value = "Pet Shop"

# Row 42
$wsOffline.Range("A42").Value = 37
$wsOffline.Range("B42").Value = "Assosiação Bichos Gerais"
$wsOffline.Range("C42").Value = "R. Pitangui, 3556 - Horto"
$wsOffline.Range("D42").Value = "(31) 3481-1968"
$wsOffline.Range("F42").Value = "Organização não Governamental"

# Row 43
$wsOffline.Range("A43").Value = 38
$wsOffline.Range("B43").Value = "Território dos Bichos"
$wsOffline.Range("C43").Value = "R. Turmalina, 566 - Prado"
$wsOffline.Range("D43").Value = "(31) 3372-2597"
$wsOffline.Range("F43").Value = "Veterinário"

# ---------------------------------------------------------------------
# Cosmetic view tweaks — widen a couple of columns on the Online sheet
# and move the scroll/selection on both sheets to where editing left off
# ---------------------------------------------------------------------

$wsOnline.Columns.Item(2).ColumnWidth = 32.833333333333336
$wsOnline.Columns.Item(5).ColumnWidth = 40.5

$wsOffline.Range("A44").Select()
$wsOnline.Select()
$wsOnline.Range("D49").Select()
